# API: Gameweeks import (#25)
# Adds two new columns ("Show Statistics Continuously" and "Gameweek") to the
# "Challenges" import-template sheet, with sample values "true" / 1 in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# New header cells
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# New sample-row cells.
# "true" must be stored as literal text (not an Excel boolean), so prefix
# with an apostrophe to force text entry, then reset the style so no
# quote-prefix formatting sticks to the cell.
$ws.Range("S2").Value = "'true"
$ws.Range("S2").Style = "Normal"
$ws.Range("T2").Value = 1
